$d = $word.ActiveDocument

$d.Content.Find.Execute("159÷6=26, 3", $true, $true, $false, $false, $false, $true, 1, $false, "832÷3=277, 1", 2) | Out-Null
$d.Content.Find.Execute("498÷7=71, 1", $true, $true, $false, $false, $false, $true, 1, $false, "469÷2=234, 1", 2) | Out-Null
$d.Content.Find.Execute("150÷9=16, 6", $true, $true, $false, $false, $false, $true, 1, $false, "940÷8=117, 4", 2) | Out-Null
$d.Content.Find.Execute("823÷7=117, 4", $true, $true, $false, $false, $false, $true, 1, $false, "977÷2=488, 1", 2) | Out-Null
$d.Content.Find.Execute("923÷5=184, 3", $true, $true, $false, $false, $false, $true, 1, $false, "178÷8=22, 2", 2) | Out-Null
$d.Content.Find.Execute("260÷2=130, 0", $true, $true, $false, $false, $false, $true, 1, $false, "100÷6=16, 4", 2) | Out-Null
$d.Content.Find.Execute("246÷7=35, 1", $true, $true, $false, $false, $false, $true, 1, $false, "827÷5=165, 2", 2) | Out-Null
$d.Content.Find.Execute("651÷7=93, 0", $true, $true, $false, $false, $false, $true, 1, $false, "513÷9=57, 0", 2) | Out-Null
$d.Content.Find.Execute("724÷4=181, 0", $true, $true, $false, $false, $false, $true, 1, $false, "646÷7=92, 2", 2) | Out-Null
$d.Content.Find.Execute("256÷4=64, 0", $true, $true, $false, $false, $false, $true, 1, $false, "923÷3=307, 2", 2) | Out-Null
$d.Content.Find.Execute("116÷6=19, 2", $true, $true, $false, $false, $false, $true, 1, $false, "186÷5=37, 1", 2) | Out-Null
$d.Content.Find.Execute("980÷3=326, 2", $true, $true, $false, $false, $false, $true, 1, $false, "990÷6=165, 0", 2) | Out-Null
$d.Content.Find.Execute("892÷5=178, 2", $true, $true, $false, $false, $false, $true, 1, $false, "884÷4=221, 0", 2) | Out-Null
$d.Content.Find.Execute("668÷4=167, 0", $true, $true, $false, $false, $false, $true, 1, $false, "530÷5=106, 0", 2) | Out-Null
$d.Content.Find.Execute("746÷7=106, 4", $true, $true, $false, $false, $false, $true, 1, $false, "986÷9=109, 5", 2) | Out-Null
$d.Content.Find.Execute("843÷4=210, 3", $true, $true, $false, $false, $false, $true, 1, $false, "926÷6=154, 2", 2) | Out-Null
$d.Content.Find.Execute("825÷7=117, 6", $true, $true, $false, $false, $false, $true, 1, $false, "526÷2=263, 0", 2) | Out-Null
$d.Content.Find.Execute("360÷4=90, 0", $true, $true, $false, $false, $false, $true, 1, $false, "523÷2=261, 1", 2) | Out-Null
$d.Content.Find.Execute("661÷4=165, 1", $true, $true, $false, $false, $false, $true, 1, $false, "844÷8=105, 4", 2) | Out-Null
$d.Content.Find.Execute("526÷5=105, 1", $true, $true, $false, $false, $false, $true, 1, $false, "223÷2=111, 1", 2) | Out-Null
$d.Content.Find.Execute("791÷8=98, 7", $true, $true, $false, $false, $false, $true, 1, $false, "133÷5=26, 3", 2) | Out-Null
$d.Content.Find.Execute("278÷2=139, 0", $true, $true, $false, $false, $false, $true, 1, $false, "382÷8=47, 6", 2) | Out-Null
$d.Content.Find.Execute("262÷6=43, 4", $true, $true, $false, $false, $false, $true, 1, $false, "135÷4=33, 3", 2) | Out-Null
$d.Content.Find.Execute("439÷7=62, 5", $true, $true, $false, $false, $false, $true, 1, $false, "608÷5=121, 3", 2) | Out-Null
$d.Content.Find.Execute("107÷6=17, 5", $true, $true, $false, $false, $false, $true, 1, $false, "995÷3=331, 2", 2) | Out-Null
